$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append the new daily data row (row 93) to the bottom of the table on Sheet1,
# mirroring the existing rows: date/day-of-week as text, hour/ranking as numbers.
$row = 93

$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/12"
$cellA.Style = "Normal"

$cellB = $ws.Cells.Item($row, 2)
$cellB.NumberFormat = "@"
$cellB.Value = "日"
$cellB.Style = "Normal"

$ws.Cells.Item($row, 3).Value = 16
$ws.Cells.Item($row, 4).Value = 36
